# Update NATMI ligand-receptor (Cxcl13-Ccr10) TPM output with refreshed
# values and rename the "Resolving-Mac" / "Neutrophils" clusters to their
# new names "MuSCs" / "ECs", expanding the sheet to the full set of
# sending/target cluster pairs (FAPs, MuSCs) x (ECs, MuSCs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 : FAPs -> ECs
$ws.Cells.Item(2, 1).Value2 = "FAPs"
$ws.Cells.Item(2, 2).Value2 = "Cxcl13"
$ws.Cells.Item(2, 3).Value2 = "Ccr10"
$ws.Cells.Item(2, 4).Value2 = "ECs"
$ws.Cells.Item(2, 5).Value2 = 3
$ws.Cells.Item(2, 6).Value2 = 1
$ws.Cells.Item(2, 7).Value2 = 7.767740666666666
$ws.Cells.Item(2, 8).Value2 = 23.303222
$ws.Cells.Item(2, 9).Value2 = 0.957755623847744
$ws.Cells.Item(2, 10).Value2 = 0.9714348434930534
$ws.Cells.Item(2, 11).Value2 = 1
$ws.Cells.Item(2, 12).Value2 = 0.5
$ws.Cells.Item(2, 13).Value2 = 0.2651005
$ws.Cells.Item(2, 14).Value2 = 0.530201
$ws.Cells.Item(2, 15).Value2 = 0.2682690254597046
$ws.Cells.Item(2, 16).Value2 = 0.2682690254597046
$ws.Cells.Item(2, 17).Value2 = 2.059231934603666
$ws.Cells.Item(2, 18).Value2 = 12.355391607622
$ws.Cells.Item(2, 19).Value2 = 0.2569361678381857
$ws.Cells.Item(2, 20).Value2 = 0.2606058787614821

# Row 3 : FAPs -> MuSCs
$ws.Cells.Item(3, 1).Value2 = "FAPs"
$ws.Cells.Item(3, 2).Value2 = "Cxcl13"
$ws.Cells.Item(3, 3).Value2 = "Ccr10"
$ws.Cells.Item(3, 4).Value2 = "MuSCs"
$ws.Cells.Item(3, 5).Value2 = 3
$ws.Cells.Item(3, 6).Value2 = 1
$ws.Cells.Item(3, 7).Value2 = 7.767740666666666
$ws.Cells.Item(3, 8).Value2 = 23.303222
$ws.Cells.Item(3, 9).Value2 = 0.957755623847744
$ws.Cells.Item(3, 10).Value2 = 0.9714348434930534
$ws.Cells.Item(3, 11).Value2 = 1
$ws.Cells.Item(3, 12).Value2 = 0.5
$ws.Cells.Item(3, 13).Value2 = 0.7230885
$ws.Cells.Item(3, 14).Value2 = 1.446177
$ws.Cells.Item(3, 15).Value2 = 0.7317309745402955
$ws.Cells.Item(3, 16).Value2 = 0.7317309745402955
$ws.Cells.Item(3, 17).Value2 = 5.616763947049
$ws.Cells.Item(3, 18).Value2 = 33.700583682294
$ws.Cells.Item(3, 19).Value2 = 0.7008194560095584
$ws.Cells.Item(3, 20).Value2 = 0.7108289647315714

# Row 4 : MuSCs -> ECs (new row)
$ws.Cells.Item(4, 1).Value2 = "MuSCs"
$ws.Cells.Item(4, 2).Value2 = "Cxcl13"
$ws.Cells.Item(4, 3).Value2 = "Ccr10"
$ws.Cells.Item(4, 4).Value2 = "ECs"
$ws.Cells.Item(4, 5).Value2 = 1
$ws.Cells.Item(4, 6).Value2 = 0.5
$ws.Cells.Item(4, 7).Value2 = 0.342617
$ws.Cells.Item(4, 8).Value2 = 0.685234
$ws.Cells.Item(4, 9).Value2 = 0.04224437615225601
$ws.Cells.Item(4, 10).Value2 = 0.02856515650694651
$ws.Cells.Item(4, 11).Value2 = 1
$ws.Cells.Item(4, 12).Value2 = 0.5
$ws.Cells.Item(4, 13).Value2 = 0.2651005
$ws.Cells.Item(4, 14).Value2 = 0.530201
$ws.Cells.Item(4, 15).Value2 = 0.2682690254597046
$ws.Cells.Item(4, 16).Value2 = 0.2682690254597046
$ws.Cells.Item(4, 17).Value2 = 0.0908279380085
$ws.Cells.Item(4, 18).Value2 = 0.363311752034
$ws.Cells.Item(4, 19).Value2 = 0.01133285762151891
$ws.Cells.Item(4, 20).Value2 = 0.00766314669822248

# Row 5 : MuSCs -> MuSCs (new row)
$ws.Cells.Item(5, 1).Value2 = "MuSCs"
$ws.Cells.Item(5, 2).Value2 = "Cxcl13"
$ws.Cells.Item(5, 3).Value2 = "Ccr10"
$ws.Cells.Item(5, 4).Value2 = "MuSCs"
$ws.Cells.Item(5, 5).Value2 = 1
$ws.Cells.Item(5, 6).Value2 = 0.5
$ws.Cells.Item(5, 7).Value2 = 0.342617
$ws.Cells.Item(5, 8).Value2 = 0.685234
$ws.Cells.Item(5, 9).Value2 = 0.04224437615225601
$ws.Cells.Item(5, 10).Value2 = 0.02856515650694651
$ws.Cells.Item(5, 11).Value2 = 1
$ws.Cells.Item(5, 12).Value2 = 0.5
$ws.Cells.Item(5, 13).Value2 = 0.7230885
$ws.Cells.Item(5, 14).Value2 = 1.446177
$ws.Cells.Item(5, 15).Value2 = 0.7317309745402955
$ws.Cells.Item(5, 16).Value2 = 0.7317309745402955
$ws.Cells.Item(5, 17).Value2 = 0.2477424126045
$ws.Cells.Item(5, 18).Value2 = 0.990969650418
$ws.Cells.Item(5, 19).Value2 = 0.03091151853073711
$ws.Cells.Item(5, 20).Value2 = 0.02090200980872403
